# Insert a new data row at row 310 (pushing existing rows 310:345 down to
# 311:346), then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("310:310").Insert()

$ws.Range("A310").Value2 = 7
$ws.Range("B310").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C310").Value2 = "Ñuble"
$ws.Range("D310").Value2 = 45124
$ws.Range("E310").Value2 = 16
$ws.Range("F310").Value2 = "Fruta"
$ws.Range("G310").Value2 = 100108
$ws.Range("H310").Value2 = "Tropicales y subtropicales"
$ws.Range("I310").Value2 = 100108005
$ws.Range("J310").Value2 = "Piña"
$ws.Range("K310").Value2 = "Caramelo"
$ws.Range("L310").Value2 = "Primera"
$ws.Range("M310").Value2 = 30
$ws.Range("N310").Value2 = 23000
$ws.Range("O310").Value2 = 23000
$ws.Range("P310").Value2 = 23000
$ws.Range("Q310").Value2 = "`$/caja 12 unidades"
$ws.Range("R310").Value2 = "Ecuador"
$ws.Range("S310").Value2 = 1917
$ws.Range("T310").Value2 = 12
